$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Neutral"
$ws.Range("B5").Value = "Positiva"
$ws.Range("B6").Value = "Negativa"
$ws.Range("B7").Value = "Neutral"
$ws.Range("B10").Value = "Neutral"
$ws.Range("B11").Value = "Positiva"
$ws.Range("B12").Value = "Negativa"
$ws.Range("B13").Value = "Neutral"
$ws.Range("B16").Value = "Neutral"
$ws.Range("B17").Value = "Positiva"
$ws.Range("B18").Value = "Negativa"
$ws.Range("B19").Value = "Neutral"
$ws.Range("B22").Value = "Negativa"
$ws.Range("B23").Value = "Neutral"
$ws.Range("B24").Value = "Positiva"
$ws.Range("B25").Value = "Positiva"
$ws.Range("B26").Value = "Positiva"
$ws.Range("B27").Value = "Negativa"
$ws.Range("B28").Value = "Negativa"
$ws.Range("B29").Value = "Positiva"
$ws.Range("B30").Value = "Neutral"
$ws.Range("B31").Value = "Positiva"
$ws.Range("B32").Value = "Positiva"
$ws.Range("B33").Value = "Negativa"
$ws.Range("B34").Value = "Positiva"
$ws.Range("B36").Value = "Negativa"
$ws.Range("B37").Value = "Neutral"
$ws.Range("B38").Value = "Positiva"
$ws.Range("B39").Value = "Negativa"
$ws.Range("B41").Value = "Negativa"
$ws.Range("B43").Value = "Neutral"
$ws.Range("B44").Value = "Positiva"
$ws.Range("B45").Value = "Positiva"
$ws.Range("B46").Value = "Positiva"
$ws.Range("B47").Value = "Neutral"
$ws.Range("B48").Value = "Positiva"
$ws.Range("B49").Value = "Positiva"
$ws.Range("B50").Value = "Neutral"
$ws.Range("B52").Value = "Positiva"
$ws.Range("B53").Value = "Neutral"
$ws.Range("B54").Value = "Positiva"
$ws.Range("B55").Value = "Neutral"
$ws.Range("B56").Value = "Positiva"
$ws.Range("B57").Value = "Neutral"
$ws.Range("B59").Value = "Neutral"
$ws.Range("B60").Value = "Neutral"
$ws.Range("B62").Value = "Neutral"
$ws.Range("B63").Value = "Positiva"
$ws.Range("B66").Value = "Neutral"
$ws.Range("B67").Value = "Neutral"
$ws.Range("B68").Value = "Neutral"
$ws.Range("B69").Value = "Negativa"
$ws.Range("B70").Value = "Negativa"
$ws.Range("B72").Value = "Positiva"
$ws.Range("B75").Value = "Positiva"
$ws.Range("B77").Value = "Positiva"
$ws.Range("B79").Value = "Positiva"
$ws.Range("B81").Value = "Positiva"
$ws.Range("B83").Value = "Positiva"
$ws.Range("B86").Value = "Positiva"
$ws.Range("B88").Value = "Positiva"
$ws.Range("B89").Value = "Positiva"
$ws.Range("B90").Value = "Negativa"
$ws.Range("B91").Value = "Positiva"
$ws.Range("B92").Value = "Negativa"
$ws.Range("B94").Value = "Negativa"
$ws.Range("B99").Value = "Neutral"
$ws.Range("B100").Value = "Negativa"
$ws.Range("B103").Value = "Negativa"
$ws.Range("B105").Value = "Negativa"
$ws.Range("B106").Value = "Positiva"
$ws.Range("B108").Value = "Neutral"
$ws.Range("B109").Value = "Neutral"
$ws.Range("B110").Value = "Neutral"
$ws.Range("B111").Value = "Positiva"
$ws.Range("B115").Value = "Negativa"
$ws.Range("B117").Value = "Positiva"
$ws.Range("B118").Value = "Positiva"
$ws.Range("B119").Value = "Negativa"
$ws.Range("B121").Value = "Positiva"
$ws.Range("B124").Value = "Negativa"
$ws.Range("B125").Value = "Negativa"
$ws.Range("B126").Value = "Neutral"
$ws.Range("B127").Value = "Positiva"
$ws.Range("B128").Value = "Negativa"
$ws.Range("B130").Value = "Positiva"
$ws.Range("B131").Value = "Neutral"
$ws.Range("B132").Value = "Positiva"
$ws.Range("B136").Value = "Neutral"
$ws.Range("B138").Value = "Neutral"
$ws.Range("B140").Value = "Negativa"
$ws.Range("B141").Value = "Positiva"
$ws.Range("B142").Value = "Negativa"
$ws.Range("B143").Value = "Positiva"
$ws.Range("B144").Value = "Negativa"
$ws.Range("B146").Value = "Negativa"
$ws.Range("B147").Value = "Positiva"
$ws.Range("B153").Value = "Positiva"
$ws.Range("B154").Value = "Positiva"
$ws.Range("B156").Value = "Negativa"
$ws.Range("B157").Value = "Positiva"
$ws.Range("B158").Value = "Neutral"
$ws.Range("B160").Value = "Positiva"
$ws.Range("B164").Value = "Positiva"
$ws.Range("B165").Value = "Negativa"
